# Fruta / hortaliza, semanal
#
# Weekly refresh of the "Cebollín baby" price series: a new week's
# observation is inserted as row 11 (pushing the existing history down by
# one row), duplicating the market/category metadata of the row that used
# to sit at row 11 (old row 10) and carrying over its volume/price figures,
# while the date is advanced to the new reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 11..46 down to 12..47, creating a blank row 11.
$ws.Rows(11).Insert()

# Seed the new row 11 with the same data as row 10 (all columns A:R), since
# that is the record being "repeated" for the new week.
for ($col = 1; $col -le 18; $col++) {
    $srcCell = $ws.Cells.Item(10, $col)
    $dstCell = $ws.Cells.Item(11, $col)
    $dstCell.Value = $srcCell.Value2()
}

# Update the new week's reporting date (column D).
$ws.Cells.Item(11, 4).Value = 44459
